$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number by Excel;
# force them to remain plain text so the literal string (e.g. trailing zeros) is kept.
$textForceRefs = @("D5", "D6", "D8", "D10", "D11", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D30", "D31", "D32", "D36", "D37", "D39", "D40", "D42", "D43", "D44", "D46", "D48", "D51")
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.599.49"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "1.665.68"
$ws.Range("E3").Value = "  -3.59%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "214.92"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").Value = "0.511"
$ws.Range("E6").Value = "  -2.35%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "23.64"
$ws.Range("E8").Value = "  -2.04%  "
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Value = "0.0622"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").Value = "0.0879"
$ws.Range("E11").Value = "  -2.30%  "
$ws.Range("D12").Value = "1.901.47"
$ws.Range("E12").Value = "  -3.51%  "
$ws.Range("D13").Value = "1.682.76"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").Value = "0.563"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "66.19"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "27.599.29"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "242.55"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "0.0₃0729"
$ws.Range("E19").Value = "  -3.60%  "
$ws.Range("D20").Value = "7.58"
$ws.Range("E20").Value = "  -4.03%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "4.48"
$ws.Range("E22").Value = "  -3.43%  "
$ws.Range("D23").Value = "9.28"
$ws.Range("E23").Value = "  -4.59%  "
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -4.58%  "
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").Value = "7.20"
$ws.Range("E26").Value = "  -4.26%  "
$ws.Range("D27").Value = "16.41"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("D30").Value = "1.23"
$ws.Range("E30").Value = "  +2.94%  "
$ws.Range("D31").Value = "0.0501"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("D33").Value = "1.469.83"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("E34").Value = "  -5.03%  "
$ws.Range("E35").Value = "  -5.89%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.37"
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "0.930"
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").Value = "0.573"
$ws.Range("E39").Value = "  -6.02%  "
$ws.Range("D40").Value = "69.37"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.22"
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.41"
$ws.Range("E44").Value = "  -7.34%  "
$ws.Range("D45").Value = "1.808.87"
$ws.Range("E45").Value = "  -3.50%  "
$ws.Range("D46").Value = "0.786"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("E47").Value = "  -4.76%  "
$ws.Range("D48").Value = "89.40"
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  -4.42%  "
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").Value = "7.89"
$ws.Range("E51").Value = "  -4.42%  "

# Restore default (General) style/number format on the forced-text cells so the
# saved workbook does not retain a lingering explicit text format on them.
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).Style = "Normal"
}
